$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing B..H to C..I
$ws.Columns("B:B").Insert()

# New column width for column B (yields stored width of 15 in OOXML)
$ws.Columns("B:B").ColumnWidth = 14.2

# New header for the inserted column
$ws.Range("B1").Value = "12-10m"

# New values for the inserted column (unstyled numbers)
$ws.Range("B2").Value = 95
$ws.Range("B3").Value = 70
$ws.Range("B4").Value = 94
$ws.Range("B5").Value = 50
$ws.Range("B6").Value = 60

# Update selection to match target
$ws.Range("D2").Select()
